$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.336
$ws.Range("A3").Value = -21.557
$ws.Range("B5").Value = 6.528
$ws.Range("A14").Value = -20.891
$ws.Range("A16").Value = -21.085
$ws.Range("B16").Value = 6.194
$ws.Range("A21").Value = -20.993
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.078
